$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACE_landing_page_data")

$ws.Range("B2").Value2 = 662.85092446876183
$ws.Range("C2").Value2 = 7937423242.1468859
$ws.Range("D2").Value2 = 11974673.262330122
$ws.Range("E2").Value2 = 0.59923449274665208
$ws.Range("F2").Value2 = 125.21093692568316
$ws.Range("G2").Value2 = 453.89943979743367
$ws.Range("H2").Value2 = -0.25232573453784612
$ws.Range("I2").Value2 = [double]"-4.8988431857763137E-2"
$ws.Range("J2").Value2 = 0.27195974513633403
$ws.Range("K2").Value2 = 0.24800162609219956
$ws.Range("L2").Value2 = [double]"-8.4930924006279773E-2"
$ws.Range("M2").Value2 = -0.24548207898660512
$ws.Range("N2").Value2 = 94.735286406577842
$ws.Range("O2").Value2 = 61.694784529334171
$ws.Range("B3").Value2 = 886.55040716031488
$ws.Range("C3").Value2 = 8346295153.5409031
$ws.Range("D3").Value2 = 9414349.2418831438
$ws.Range("E3").Value2 = 0.48015521792467764
$ws.Range("F3").Value2 = 136.83222415718748
$ws.Range("G3").Value2 = 601.57542605190906
$ws.Range("H3").Value2 = 1.2156413238068402
$ws.Range("I3").Value2 = [double]"-4.0950957368767971E-2"
$ws.Range("J3").Value2 = -0.56714607534786987
$ws.Range("K3").Value2 = -0.5071695245565484
$ws.Range("L3").Value2 = [double]"8.605977981439894E-2"
$ws.Range("M3").Value2 = 1.2213343463924975
$ws.Range("N3").Value2 = 99.615282905169551
$ws.Range("O3").Value2 = 48.503724088156154
$ws.Range("B4").Value2 = 400.13263773085521
$ws.Range("D4").Value2 = 21749483.383913253
$ws.Range("E4").Value2 = 0.97428069457886368
$ws.Range("F4").Value2 = 125.98958795856642
$ws.Range("G4").Value2 = 270.81714512220213
$ws.Range("H4").Value2 = [double]"-2.4712782600395666E-3"
$ws.Range("I4").Value2 = [double]"1.408073099924545E-2"
$ws.Range("J4").Value2 = [double]"1.659301521705947E-2"
$ws.Range("K4").Value2 = [double]"1.017151917622483E-2"
$ws.Range("L4").Value2 = [double]"6.3459591994670106E-3"
$ws.Range("M4").Value2 = [double]"-1.8417737172722326E-3"
$ws.Range("N4").Value2 = 103.86881012035276
$ws.Range("O4").Value2 = 112.05564123540506
$ws.Range("B5").Value2 = 401.12392657016977
$ws.Range("C5").Value2 = 8581839581.0689077
$ws.Range("D5").Value2 = 21394484.379050527
$ws.Range("E5").Value2 = 0.96447056374482887
$ws.Range("F5").Value2 = 125.19510493070318
$ws.Range("G5").Value2 = 271.31684936441462
$ws.Range("H5").Value2 = [double]"-3.6096055075297007E-2"
$ws.Range("J5").Value2 = [double]"5.3266327393592938E-2"
$ws.Range("K5").Value2 = [double]"4.9170449720506104E-2"
$ws.Range("L5").Value2 = [double]"-4.4919207430083041E-4"
$ws.Range("M5").Value2 = [double]"-3.0644875722143139E-2"
$ws.Range("N5").Value2 = 102.42656915293469
$ws.Range("O5").Value2 = 110.22664877495674
$ws.Range("B6").Value2 = 416.14512388109819
$ws.Range("D6").Value2 = 20312511.491745111
$ws.Range("E6").Value2 = 0.9192696610944
$ws.Range("F6").Value2 = 125.25136685198844
$ws.Range("G6").Value2 = 279.89417146429003
$ws.Range("H6").Value2 = [double]"-3.5966192083694559E-2"
$ws.Range("J6").Value2 = [double]"4.652209899962334E-2"
$ws.Range("K6").Value2 = [double]"4.2160412093241728E-2"
$ws.Range("L6").Value2 = [double]"1.1290411509541132E-2"
$ws.Range("M6").Value2 = [double]"-3.9024986400600481E-2"
$ws.Range("O6").Value2 = 104.65220989996233
$ws.Range("B7").Value2 = 431.67067426874587
$ws.Range("C7").Value2 = 8378528786.0762291
$ws.Range("D7").Value2 = 19409538.996990088
$ws.Range("E7").Value2 = 0.88208077223734849
$ws.Range("F7").Value2 = 123.85301534207885
$ws.Range("G7").Value2 = 291.26061292262608
$ws.Range("H7").Value2 = [double]"-1.9201466363056863E-2"
$ws.Range("J7").Value2 = [double]"2.6274287320018574E-2"
$ws.Range("K7").Value2 = [double]"1.8152756582704876E-2"
$ws.Range("L7").Value2 = [double]"1.6252900049103625E-2"
$ws.Range("M7").Value2 = [double]"-2.7345176258441239E-2"
